$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell as genuine text (not auto-converted to a
# number), without adding any new cell style (keeps output matching the
# "numberStoredAsText" convention already used throughout this workbook).
function Set-TextValue {
    param($cell, [string]$text)

    # Escape embedded double quotes the way Excel formulas require ("" for a literal ").
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null  # xlPasteValues
}

# --- Orders sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("Orders")

# Row 51 gains a Number value that was previously missing.
Set-TextValue $ws.Cells.Item(51, 6) '10'

# New row 52
Set-TextValue $ws.Cells.Item(52, 3) '77_珍爱mini_undefined_Gerbera L._20stems'
Set-TextValue $ws.Cells.Item(52, 6) '8'

# New row 53
Set-TextValue $ws.Cells.Item(53, 3) '846_玛格丽特_undefined_undefined_1bunch'
Set-TextValue $ws.Cells.Item(53, 6) '10'

# New row 54
Set-TextValue $ws.Cells.Item(54, 3) '846_玛格丽特_undefined_undefined_1bunch'
Set-TextValue $ws.Cells.Item(54, 6) '5'

# --- Summary sheet ------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")
Set-TextValue $ws2.Cells.Item(2, 7) '014613710139878631015151015510105101015551075555510101515101010555581051010108105'
